$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: quantity 5 -> 6, amount 22.45 -> 26.94
$ws.Range("K19").Value = 6
$ws.Range("L19").Value = 26.94

# Row 25: quantity 9 -> 10, amount 15.3 -> 16.8
$ws.Range("K25").Value = 10
$ws.Range("L25").Value = 16.8

# Row 35: totals 379 -> 381, 486.85 -> 492.84
$ws.Range("K35").Value = 381
$ws.Range("L35").Value = 492.84
